$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 919.2778  # ALC!H28
$ws.Cells.Item(28, 9).Value = 368.125  # ALC!I28
$ws.Cells.Item(28, 10).Value = 1360.2  # ALC!J28
$ws.Cells.Item(28, 11).Value = 368.125  # ALC!K28
$ws.Cells.Item(28, 12).Value = 1360.2  # ALC!L28
$ws.Cells.Item(28, 13).Value = 116.875  # ALC!M28
$ws.Cells.Item(28, 14).Value = -2330.2  # ALC!N28

$ws.Cells.Item(43, 8).Value = 935.3333  # ALC!H43
$ws.Cells.Item(43, 9).Value = 1000  # ALC!I43
$ws.Cells.Item(43, 10).Value = 922.4  # ALC!J43
$ws.Cells.Item(43, 11).Value = 1000  # ALC!K43
$ws.Cells.Item(43, 12).Value = 922.4  # ALC!L43
$ws.Cells.Item(43, 13).Value = -931  # ALC!M43
$ws.Cells.Item(43, 14).Value = -1060.4  # ALC!N43

$ws.Cells.Item(111, 8).Value = 2981.3635  # ALC!H111
$ws.Cells.Item(111, 9).Value = 3607.5  # ALC!I111
$ws.Cells.Item(111, 10).Value = 2230  # ALC!J111
$ws.Cells.Item(111, 11).Value = 10822.5  # ALC!K111
$ws.Cells.Item(111, 12).Value = 6690  # ALC!L111
$ws.Cells.Item(111, 13).Value = -7755.5  # ALC!M111
$ws.Cells.Item(111, 14).Value = -12824  # ALC!N111

$ws.Cells.Item(112, 8).Value = 4630628.5  # ALC!H112
$ws.Cells.Item(112, 10).Value = 5849033.5  # ALC!J112
$ws.Cells.Item(112, 12).Value = 17547100.5  # ALC!L112
$ws.Cells.Item(112, 14).Value = -17549316.5  # ALC!N112

$ws.Cells.Item(132, 8).Value = 2301.279  # ALC!H132
$ws.Cells.Item(132, 9).Value = 2393.878  # ALC!I132
$ws.Cells.Item(132, 10).Value = 403  # ALC!J132
$ws.Cells.Item(132, 11).Value = 7181.634  # ALC!K132
$ws.Cells.Item(132, 12).Value = 1209  # ALC!L132
$ws.Cells.Item(132, 13).Value = -4651.634  # ALC!M132
$ws.Cells.Item(132, 14).Value = -6269  # ALC!N132

$ws.Cells.Item(134, 8).Value = 46994.5  # ALC!H134
$ws.Cells.Item(134, 10).Value = 46994.5  # ALC!J134
$ws.Cells.Item(134, 12).Value = 46994.5  # ALC!L134
$ws.Cells.Item(134, 14).Value = -57134.5  # ALC!N134

$ws.Cells.Item(135, 8).Value = 17862888  # ALC!H135
$ws.Cells.Item(135, 9).Value = 765.64703  # ALC!I135
$ws.Cells.Item(135, 11).Value = 6890.82327  # ALC!K135
$ws.Cells.Item(135, 13).Value = -4355.82327  # ALC!M135

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5923.7856  # ARM!H32
$ws.Cells.Item(32, 9).Value = 5090.731  # ARM!I32
$ws.Cells.Item(32, 11).Value = 5090.731  # ARM!K32
$ws.Cells.Item(32, 13).Value = -4803.731  # ARM!M32

$ws.Cells.Item(45, 8).Value = 2327.842  # ARM!H45
$ws.Cells.Item(45, 10).Value = 2754.2  # ARM!J45
$ws.Cells.Item(45, 12).Value = 2754.2  # ARM!L45
$ws.Cells.Item(45, 14).Value = -3508.2  # ARM!N45

$ws.Cells.Item(61, 8).Value = 2497.8386  # ARM!H61
$ws.Cells.Item(61, 9).Value = 2733.96  # ARM!I61
$ws.Cells.Item(61, 11).Value = 2733.96  # ARM!K61
$ws.Cells.Item(61, 13).Value = -2521.96  # ARM!M61

$ws.Cells.Item(97, 8).Value = 62500790  # ARM!H97
$ws.Cells.Item(97, 9).Value = 611.6923  # ARM!I97
$ws.Cells.Item(97, 11).Value = 611.6923  # ARM!K97
$ws.Cells.Item(97, 13).Value = -115.6923  # ARM!M97

$ws.Cells.Item(110, 8).Value = 809.6  # ARM!H110
$ws.Cells.Item(110, 9).Value = 730.8570999999999  # ARM!I110
$ws.Cells.Item(110, 10).Value = 993.3333  # ARM!J110
$ws.Cells.Item(110, 11).Value = 730.8570999999999  # ARM!K110
$ws.Cells.Item(110, 12).Value = 993.3333  # ARM!L110
$ws.Cells.Item(110, 13).Value = 1314.1429  # ARM!M110
$ws.Cells.Item(110, 14).Value = -5083.3333  # ARM!N110

$ws.Cells.Item(123, 8).Value = 0  # ARM!H123
$ws.Cells.Item(123, 10).Value = 0  # ARM!J123
$ws.Cells.Item(123, 14).ClearContents()  # ARM!N123

$ws.Cells.Item(136, 8).Value = 2497.8386  # ARM!H136
$ws.Cells.Item(136, 9).Value = 2733.96  # ARM!I136
$ws.Cells.Item(136, 11).Value = 8201.880000000001  # ARM!K136
$ws.Cells.Item(136, 13).Value = -5651.880000000001  # ARM!M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1102.5  # BSM!H107
$ws.Cells.Item(107, 9).Value = 1205.5  # BSM!I107
$ws.Cells.Item(107, 10).Value = 999.5  # BSM!J107
$ws.Cells.Item(107, 11).Value = 1205.5  # BSM!K107
$ws.Cells.Item(107, 12).Value = 999.5  # BSM!L107
$ws.Cells.Item(107, 13).Value = 714.5  # BSM!M107
$ws.Cells.Item(107, 14).Value = -4839.5  # BSM!N107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3364.7778  # CRP!H31
$ws.Cells.Item(31, 9).Value = 2306.4736  # CRP!I31
$ws.Cells.Item(31, 10).Value = 4547.5884  # CRP!J31
$ws.Cells.Item(31, 11).Value = 2306.4736  # CRP!K31
$ws.Cells.Item(31, 12).Value = 4547.5884  # CRP!L31
$ws.Cells.Item(31, 13).Value = -2011.4736  # CRP!M31
$ws.Cells.Item(31, 14).Value = -5137.5884  # CRP!N31

$ws.Cells.Item(34, 8).Value = 3364.7778  # CRP!H34
$ws.Cells.Item(34, 9).Value = 2306.4736  # CRP!I34
$ws.Cells.Item(34, 10).Value = 4547.5884  # CRP!J34
$ws.Cells.Item(34, 11).Value = 2306.4736  # CRP!K34
$ws.Cells.Item(34, 12).Value = 4547.5884  # CRP!L34
$ws.Cells.Item(34, 13).Value = -2104.4736  # CRP!M34
$ws.Cells.Item(34, 14).Value = -4951.5884  # CRP!N34

$ws.Cells.Item(41, 8).Value = 20000  # CRP!H41
$ws.Cells.Item(41, 10).Value = 20000  # CRP!J41
$ws.Cells.Item(41, 12).Value = 20000  # CRP!L41
$ws.Cells.Item(41, 14).Value = -20856  # CRP!N41

$ws.Cells.Item(50, 8).Value = 16533.334  # CRP!H50
$ws.Cells.Item(50, 10).Value = 16533.334  # CRP!J50
$ws.Cells.Item(50, 12).Value = 16533.334  # CRP!L50
$ws.Cells.Item(50, 14).Value = -17783.334  # CRP!N50

$ws.Cells.Item(51, 8).Value = 20000  # CRP!H51
$ws.Cells.Item(51, 10).Value = 20000  # CRP!J51
$ws.Cells.Item(51, 12).Value = 20000  # CRP!L51
$ws.Cells.Item(51, 14).Value = -21472  # CRP!N51

$ws.Cells.Item(60, 8).Value = 6312  # CRP!H60
$ws.Cells.Item(60, 9).Value = 1986.6666  # CRP!I60
$ws.Cells.Item(60, 10).Value = 12800  # CRP!J60
$ws.Cells.Item(60, 11).Value = 1986.6666  # CRP!K60
$ws.Cells.Item(60, 12).Value = 12800  # CRP!L60
$ws.Cells.Item(60, 13).Value = -1475.6666  # CRP!M60
$ws.Cells.Item(60, 14).Value = -13822  # CRP!N60

$ws.Cells.Item(61, 8).Value = 20000  # CRP!H61
$ws.Cells.Item(61, 10).Value = 20000  # CRP!J61
$ws.Cells.Item(61, 12).Value = 20000  # CRP!L61
$ws.Cells.Item(61, 14).Value = -20696  # CRP!N61

$ws.Cells.Item(62, 8).Value = 47622300  # CRP!H62
$ws.Cells.Item(62, 9).Value = 52634540  # CRP!I62
$ws.Cells.Item(62, 10).Value = 6003  # CRP!J62
$ws.Cells.Item(62, 11).Value = 52634540  # CRP!K62
$ws.Cells.Item(62, 12).Value = 6003  # CRP!L62
$ws.Cells.Item(62, 13).Value = -52633916  # CRP!M62
$ws.Cells.Item(62, 14).Value = -7251  # CRP!N62

$ws.Cells.Item(65, 8).Value = 47622300  # CRP!H65
$ws.Cells.Item(65, 9).Value = 52634540  # CRP!I65
$ws.Cells.Item(65, 10).Value = 6003  # CRP!J65
$ws.Cells.Item(65, 11).Value = 263172700  # CRP!K65
$ws.Cells.Item(65, 12).Value = 30015  # CRP!L65
$ws.Cells.Item(65, 13).Value = -263169580  # CRP!M65
$ws.Cells.Item(65, 14).Value = -36255  # CRP!N65

$ws.Cells.Item(99, 8).Value = 22730724  # CRP!H99
$ws.Cells.Item(99, 9).Value = 3190  # CRP!I99
$ws.Cells.Item(99, 11).Value = 3190  # CRP!K99
$ws.Cells.Item(99, 13).Value = -1692  # CRP!M99

$ws.Cells.Item(107, 8).Value = 1120.3928  # CRP!H107
$ws.Cells.Item(107, 9).Value = 494.5  # CRP!I107
$ws.Cells.Item(107, 11).Value = 494.5  # CRP!K107
$ws.Cells.Item(107, 13).Value = 1425.5  # CRP!M107

$ws.Cells.Item(126, 8).Value = 22730724  # CRP!H126
$ws.Cells.Item(126, 9).Value = 3190  # CRP!I126
$ws.Cells.Item(126, 11).Value = 9570  # CRP!K126
$ws.Cells.Item(126, 13).Value = -7100  # CRP!M126

$ws.Cells.Item(132, 8).Value = 4824.0713  # CRP!H132
$ws.Cells.Item(132, 9).Value = 3853.125  # CRP!I132
$ws.Cells.Item(132, 11).Value = 11559.375  # CRP!K132
$ws.Cells.Item(132, 13).Value = -9029.375  # CRP!M132

$ws.Cells.Item(134, 8).Value = 1262.4286  # CRP!H134
$ws.Cells.Item(134, 9).Value = 1219.7142  # CRP!I134
$ws.Cells.Item(134, 10).Value = 1305.1428  # CRP!J134
$ws.Cells.Item(134, 11).Value = 3659.1426  # CRP!K134
$ws.Cells.Item(134, 12).Value = 3915.4284  # CRP!L134
$ws.Cells.Item(134, 13).Value = -1124.1426  # CRP!M134
$ws.Cells.Item(134, 14).Value = -8985.428400000001  # CRP!N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 605.05  # CUL!H122
$ws.Cells.Item(122, 10).Value = 674.8125  # CUL!J122
$ws.Cells.Item(122, 12).Value = 6073.3125  # CUL!L122
$ws.Cells.Item(122, 14).Value = -10973.3125  # CUL!N122

$ws.Cells.Item(129, 8).Value = 232866.45  # CUL!H129
$ws.Cells.Item(129, 9).Value = 620.9091  # CUL!I129
$ws.Cells.Item(129, 10).Value = 465112  # CUL!J129
$ws.Cells.Item(129, 11).Value = 1862.7273  # CUL!K129
$ws.Cells.Item(129, 12).Value = 1395336  # CUL!L129
$ws.Cells.Item(129, 13).Value = 3137.2727  # CUL!M129
$ws.Cells.Item(129, 14).Value = -1405336  # CUL!N129

$ws.Cells.Item(131, 8).Value = 108301.59  # CUL!H131
$ws.Cells.Item(131, 10).Value = 119844.414  # CUL!J131
$ws.Cells.Item(131, 12).Value = 359533.242  # CUL!L131
$ws.Cells.Item(131, 14).Value = -369613.242  # CUL!N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 4315.857  # GSM!H113
$ws.Cells.Item(113, 9).Value = 2737  # GSM!I113
$ws.Cells.Item(113, 10).Value = 5500  # GSM!J113
$ws.Cells.Item(113, 11).Value = 2737  # GSM!K113
$ws.Cells.Item(113, 12).Value = 5500  # GSM!L113
$ws.Cells.Item(113, 13).Value = -567  # GSM!M113
$ws.Cells.Item(113, 14).Value = -9840  # GSM!N113

$ws.Cells.Item(132, 8).Value = 36484.266  # GSM!H132
$ws.Cells.Item(132, 9).Value = 3470.8333  # GSM!I132
$ws.Cells.Item(132, 11).Value = 10412.4999  # GSM!K132
$ws.Cells.Item(132, 13).Value = -7882.499899999999  # GSM!M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1887.25  # LTW!H136
$ws.Cells.Item(136, 9).Value = 2574.5  # LTW!I136
$ws.Cells.Item(136, 10).Value = 1200  # LTW!J136
$ws.Cells.Item(136, 11).Value = 7723.5  # LTW!K136
$ws.Cells.Item(136, 12).Value = 3600  # LTW!L136
$ws.Cells.Item(136, 13).Value = -5173.5  # LTW!M136
$ws.Cells.Item(136, 14).Value = -8700  # LTW!N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2841517  # WVR!H107
$ws.Cells.Item(107, 9).Value = 500.66666  # WVR!I107
$ws.Cells.Item(107, 10).Value = 6494252  # WVR!J107
$ws.Cells.Item(107, 11).Value = 1501.99998  # WVR!K107
$ws.Cells.Item(107, 12).Value = 19482756  # WVR!L107
$ws.Cells.Item(107, 13).Value = 418.0000199999999  # WVR!M107
$ws.Cells.Item(107, 14).Value = -19486596  # WVR!N107

$ws.Cells.Item(113, 8).Value = 3862479.2  # WVR!H113
$ws.Cells.Item(113, 9).Value = 9000  # WVR!I113
$ws.Cells.Item(113, 10).Value = 4504725.5  # WVR!J113
$ws.Cells.Item(113, 11).Value = 27000  # WVR!K113
$ws.Cells.Item(113, 12).Value = 13514176.5  # WVR!L113
$ws.Cells.Item(113, 13).Value = -24830  # WVR!M113
$ws.Cells.Item(113, 14).Value = -13518516.5  # WVR!N113

$ws.Cells.Item(136, 8).Value = 41670710  # WVR!H136
$ws.Cells.Item(136, 9).Value = 66668796  # WVR!I136
$ws.Cells.Item(136, 10).Value = 7236.6665  # WVR!J136
$ws.Cells.Item(136, 11).Value = 200006388  # WVR!K136
$ws.Cells.Item(136, 12).Value = 21709.9995  # WVR!L136
$ws.Cells.Item(136, 13).Value = -200003838  # WVR!M136
$ws.Cells.Item(136, 14).Value = -26809.9995  # WVR!N136
